$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some updated Price cells are plain decimal numbers (e.g. "570.21"). Assigning
# such a string straight to .Value lets Excel auto-detect it as a Number, but the
# source data keeps these as text cells, so mark them Text first (per-cell, since
# NumberFormat on a multi-area Range only actually applies to its first area here).
$textRefs = @("D5", "D6", "D10", "D12", "D14", "D16", "D21", "D24", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D37", "D39", "D42", "D43", "D47", "D49", "D51")
foreach ($r in $textRefs) {
    $ws.Range($r).NumberFormat = "@"
}

$ws.Range('D2').Value = '68.834.43'
$ws.Range('E2').Value = '  -1.35%  '
$ws.Range('D3').Value = '3.493.47'
$ws.Range('E3').Value = '  -2.12%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '570.21'
$ws.Range('E5').Value = '  -1.11%  '
$ws.Range('D6').Value = '182.77'
$ws.Range('E6').Value = '  -2.77%  '
$ws.Range('E7').Value = '  -2.77%  '
$ws.Range('D8').Value = '3.485.93'
$ws.Range('E8').Value = '  -2.20%  '
$ws.Range('E9').Value = '  +0.10%  '
$ws.Range('D10').Value = '0.184'
$ws.Range('E10').Value = '  +4.11%  '
$ws.Range('E11').Value = '  -2.19%  '
$ws.Range('D12').Value = '53.91'
$ws.Range('E12').Value = '  -3.17%  '
$ws.Range('E13').Value = '  -0.22%  '
$ws.Range('D14').Value = '9.41'
$ws.Range('E14').Value = '  -1.48%  '
$ws.Range('D15').Value = '4.053.94'
$ws.Range('E15').Value = '  -2.14%  '
$ws.Range('D16').Value = '19.22'
$ws.Range('E16').Value = '  -2.54%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.494.31'
$ws.Range('E17').Value = '  -2.03%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '68.776.97'
$ws.Range('E18').Value = '  -1.27%  '
$ws.Range('E19').Value = '  -2.63%  '
$ws.Range('E20').Value = '  -1.04%  '
$ws.Range('D21').Value = '540.02'
$ws.Range('E21').Value = '  +13.84%  '
$ws.Range('E22').Value = '  -2.27%  '
$ws.Range('E23').Value = '  +0.87%  '
$ws.Range('D24').Value = '4.99'
$ws.Range('E24').Value = '  -1.20%  '
$ws.Range('D25').Value = '4.39'
$ws.Range('E25').Value = '  +1.05%  '
$ws.Range('D26').Value = '94.14'
$ws.Range('E26').Value = '  -0.51%  '
$ws.Range('E27').Value = '  -3.58%  '
$ws.Range('D28').Value = '10.76'
$ws.Range('E28').Value = '  -1.56%  '
$ws.Range('D29').Value = '9.07'
$ws.Range('E29').Value = '  -2.25%  '
$ws.Range('D30').Value = '31.48'
$ws.Range('D31').Value = '7.17'
$ws.Range('E31').Value = '  -7.35%  '
$ws.Range('D32').Value = '12.51'
$ws.Range('E32').Value = '  +3.12%  '
$ws.Range('D33').Value = '64.62'
$ws.Range('E33').Value = '  -2.46%  '
$ws.Range('E34').Value = '  -4.62%  '
$ws.Range('D35').Value = '568.05'
$ws.Range('E35').Value = '  -4.22%  '
$ws.Range('E36').Value = '  +0.10%  '
$ws.Range('D37').Value = '37.70'
$ws.Range('E37').Value = '  -2.88%  '
$ws.Range('E38').Value = '  +0.46%  '
$ws.Range('D39').Value = '2.98'
$ws.Range('E39').Value = '  +5.53%  '
$ws.Range('D40').Value = '0.0₃0763'
$ws.Range('E40').Value = '  -3.96%  '
$ws.Range('E41').Value = '  -3.67%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').Value = '0.131'
$ws.Range('E42').Value = '  -4.88%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = '3.31'
$ws.Range('E43').Value = '  -4.22%  '
$ws.Range('D44').Value = '3.248.43'
$ws.Range('E44').Value = '  +0.90%  '
$ws.Range('E45').Value = '  -3.64%  '
$ws.Range('E46').Value = '  +1.99%  '
$ws.Range('D47').Value = '0.0436'
$ws.Range('E47').Value = '  -0.52%  '
$ws.Range('E48').Value = '  -2.22%  '
$ws.Range('D49').Value = '8.96'
$ws.Range('E49').Value = '  -5.42%  '
$ws.Range('E50').Value = '  -0.17%  '
$ws.Range('D51').Value = '137.48'
$ws.Range('E51').Value = '  +1.80%  '

# Drop back to the default (unstyled) look now that the values are locked in as
# text, so these cells end up with no explicit style - matching the rest of the sheet.
foreach ($r in $textRefs) {
    $ws.Range($r).Style = "Normal"
}

